$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Status column (E) updates for the audio-event asset list.
# Almost every event is now marked "DONE"; the two Vomit / Voice Line 1
# rows (still waiting on an unimplemented feature) are marked as
# "POSTPONED (Feature not in game yet)" instead of the old
# "FIXING BUGS" placeholder text.  The old "Made" status string is no
# longer used anywhere.
# ------------------------------------------------------------------
$postponedText = "POSTPONED (Feature not in game yet)"
$doneText = "DONE"

$postponedRows = @(15, 16)

for ($r = 2; $r -le 33; $r++) {
    if ($postponedRows -contains $r) {
        $ws.Range("E$r").Value = $postponedText
    } else {
        $ws.Range("E$r").Value = $doneText
    }
}

# ------------------------------------------------------------------
# View / layout tweaks
# ------------------------------------------------------------------
# Widen the Status column so the longer status text fits.
$ws.Columns("E").ColumnWidth = 36.67

# Scroll back to the top of the sheet and move the selection.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F4").Select()
